$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row values: target row -> (Date serial D, Volumen M, Precio N=O=P=S)
$rowData = @{
    2 = @(44431, 100, 1300)
    3 = @(44424, 50, 1200)
    4 = @(44760, 80, 2300)
    5 = @(44438, 60, 1200)
    6 = @(44749, 120, 2300)
    7 = @(44435, 130, 1300)
    8 = @(44432, 30, 1300)
    9 = @(44748, 300, 2300)
    10 = @(44473, 120, 1200)
    11 = @(44357, 35, 1000)
    12 = @(44405, 50, 1200)
    13 = @(44753, 160, 2300)
    14 = @(44343, 60, 1300)
    15 = @(44476, 80, 1200)
    16 = @(44763, 50, 2300)
    17 = @(44811, 60, 2500)
    18 = @(44812, 50, 2500)
    19 = @(44418, 40, 1200)
    20 = @(44830, 50, 2500)
    21 = @(44417, 80, 1200)
    22 = @(44762, 50, 2300)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $d = $vals[0]
    $m = $vals[1]
    $n = $vals[2]
    $ws.Cells.Item($r, 4).Value = $d   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $m  # M: Volumen
    $ws.Cells.Item($r, 14).Value = $n  # N: Precio mínimo
    $ws.Cells.Item($r, 15).Value = $n  # O: Precio máximo
    $ws.Cells.Item($r, 16).Value = $n  # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $n  # S: Precio $/Kg
}
